$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.529.39'
$ws.Range("E2").Value = '  +3.53%  '

$ws.Range("D3").Value = '1.823.30'
$ws.Range("E3").Value = '  +4.57%  '

$ws.Range("E4").Value = '  -0.05%  '

$ws.Range("D5").Formula = '="343.04"'
$ws.Range("D5").Copy()
$ws.Range("D5").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E5").Value = '  +2.42%  '

$ws.Range("D6").Formula = '="1.001"'
$ws.Range("D6").Copy()
$ws.Range("D6").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E6").Value = '  -0.09%  '

$ws.Range("D7").Formula = '="0.3818"'
$ws.Range("D7").Copy()
$ws.Range("D7").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E7").Value = '  +0.47%  '

$ws.Range("E8").Value = '  +4.23%  '

$ws.Range("D9").Formula = '="49.95"'
$ws.Range("D9").Copy()
$ws.Range("D9").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E9").Value = '  +2.78%  '

$ws.Range("D10").Formula = '="1.237"'
$ws.Range("D10").Copy()
$ws.Range("D10").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E10").Value = '  +3.86%  '

$ws.Range("D11").Formula = '="0.07737"'
$ws.Range("D11").Copy()
$ws.Range("D11").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E11").Value = '  +3.41%  '

$ws.Range("D12").Formula = '="1.003"'
$ws.Range("D12").Copy()
$ws.Range("D12").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E12").Value = '  -0.04%  '

$ws.Range("D13").Formula = '="22.19"'
$ws.Range("D13").Copy()
$ws.Range("D13").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E13").Value = '  +8.48%  '

$ws.Range("D14").Formula = '="6.612"'
$ws.Range("D14").Copy()
$ws.Range("D14").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E14").Value = '  +2.24%  '

$ws.Range("D15").Value = '1.825.94'
$ws.Range("E15").Value = '  +4.79%  '

$ws.Range("D16").Formula = '="7.216"'
$ws.Range("D16").Copy()
$ws.Range("D16").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E16").Value = '  +1.21%  '

$ws.Range("E17").Value = '  +3.33%  '

$ws.Range("D18").Formula = '="0.06742"'
$ws.Range("D18").Copy()
$ws.Range("D18").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E18").Value = '  +0.50%  '

$ws.Range("D19").Formula = '="86.99"'
$ws.Range("D19").Copy()
$ws.Range("D19").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E19").Value = '  +4.13%  '

$ws.Range("E20").Value = '  +0.03%  '

$ws.Range("D21").Formula = '="17.58"'
$ws.Range("D21").Copy()
$ws.Range("D21").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E21").Value = '  +4.82%  '

$ws.Range("D22").Formula = '="6.532"'
$ws.Range("D22").Copy()
$ws.Range("D22").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E22").Value = '  +5.05%  '

$ws.Range("D23").Formula = '="13.16"'
$ws.Range("D23").Copy()
$ws.Range("D23").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E23").Value = '  +0.73%  '

$ws.Range("D24").Value = '27.524.94'
$ws.Range("E24").Value = '  +3.52%  '

$ws.Range("D25").Formula = '="2.477"'
$ws.Range("D25").Copy()
$ws.Range("D25").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E25").Value = '  +1.13%  '

$ws.Range("D26").Formula = '="2.686"'
$ws.Range("D26").Copy()
$ws.Range("D26").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E26").Value = '  +8.60%  '

$ws.Range("D27").Formula = '="22.00"'
$ws.Range("D27").Copy()
$ws.Range("D27").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E27").Value = '  +11.66%  '

$ws.Range("D28").Formula = '="1.482"'
$ws.Range("D28").Copy()
$ws.Range("D28").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E28").Value = '  +4.00%  '

$ws.Range("D29").Formula = '="152.94"'
$ws.Range("D29").Copy()
$ws.Range("D29").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E29").Value = '  -0.86%  '

$ws.Range("D30").Value = '2.031.18'
$ws.Range("E30").Value = '  +4.99%  '

$ws.Range("D31").Formula = '="135.43"'
$ws.Range("D31").Copy()
$ws.Range("D31").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E31").Value = '  +2.48%  '

$ws.Range("D32").Formula = '="6.339"'
$ws.Range("D32").Copy()
$ws.Range("D32").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E32").Value = '  +3.51%  '

$ws.Range("D33").Formula = '="4.094"'
$ws.Range("D33").Copy()
$ws.Range("D33").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E33").Value = '  -1.27%  '

$ws.Range("E34").Value = '  +6.61%  '

$ws.Range("D35").Formula = '="0.08794"'
$ws.Range("D35").Copy()
$ws.Range("D35").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E35").Value = '  +0.94%  '

$ws.Range("D36").Formula = '="1.697"'
$ws.Range("D36").Copy()
$ws.Range("D36").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E36").Value = '  -0.98%  '

$ws.Range("D37").Formula = '="5.622"'
$ws.Range("D37").Copy()
$ws.Range("D37").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E37").Value = '  +3.14%  '

$ws.Range("D38").Formula = '="0.7003"'
$ws.Range("D38").Copy()
$ws.Range("D38").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E38").Value = '  +11.79%  '

$ws.Range("D39").Formula = '="9.114"'
$ws.Range("D39").Copy()
$ws.Range("D39").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E39").Value = '  +5.66%  '

$ws.Range("D40").Formula = '="0.06520"'
$ws.Range("D40").Copy()
$ws.Range("D40").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E40").Value = '  +2.83%  '

$ws.Range("D41").Formula = '="0.2260"'
$ws.Range("D41").Copy()
$ws.Range("D41").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E41").Value = '  +3.44%  '

$ws.Range("D42").Formula = '="0.02407"'
$ws.Range("D42").Copy()
$ws.Range("D42").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E42").Value = '  +1.52%  '

$ws.Range("D43").Formula = '="1.304"'
$ws.Range("D43").Copy()
$ws.Range("D43").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E43").Value = '  +5.98%  '

$ws.Range("D44").Formula = '="14.76"'
$ws.Range("D44").Copy()
$ws.Range("D44").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E44").Value = '  +2.94%  '

$ws.Range("D45").Formula = '="0.6617"'
$ws.Range("D45").Copy()
$ws.Range("D45").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E45").Value = '  +8.78%  '

$ws.Range("D46").Formula = '="1.001"'
$ws.Range("D46").Copy()
$ws.Range("D46").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E46").Value = '  -0.09%  '

$ws.Range("D47").Formula = '="3.947"'
$ws.Range("D47").Copy()
$ws.Range("D47").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E47").Value = '  +0.39%  '

$ws.Range("D48").Formula = '="2.189"'
$ws.Range("D48").Copy()
$ws.Range("D48").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E48").Value = '  +5.82%  '

$ws.Range("D49").Formula = '="133.24"'
$ws.Range("D49").Copy()
$ws.Range("D49").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E49").Value = '  +3.36%  '

$ws.Range("D50").Formula = '="0.07307"'
$ws.Range("D50").Copy()
$ws.Range("D50").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("E50").Value = '  +0.82%  '

$ws.Range("E51").Value = '  +3.86%  '

